{"js": "// Office.js (Word JavaScript API) edit script.\n//\n// Target change (per the commit \"9/21/18 - up to exercises\"):\n//   1) In the \"Use for..in loop ...\" bullet, split the single \"for..in\"\n//      run into \"for..\" / \"in\" and wrap \"for..\" in its own\n//      gramStart/gramEnd proofErr pair (matching the existing\n//      \"for..of\" / \"if..else\" bullets elsewhere in the doc), and drop\n//      the stray _GoBack bookmark that was sitting at the end of that\n//      paragraph.\n//   2) The next (previously empty) bullet in the same list gets new\n//      text about the `continue` keyword, with the _GoBack bookmark\n//      relocated into the middle of that new sentence.\n//\n// Word's JS API has no direct way to create proofErr/bookmark markup,\n// so we use Range.insertOoxml() with a literal FlatOPC-wrapped <w:p>\n// fragment for each of the two affected paragraphs; that lets us\n// control the run/proofErr/bookmark structure exactly while keeping\n// each paragraph's own pPr (style + numbering) and w14 identity\n// attributes intact.\n\nfunction flatOpcPackage(paragraphXml) {\n  return (\n    '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">' +\n    '<w:body>' + paragraphXml + '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>'\n  );\n}\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// --- Paragraph: \"Use for..in loop to iterate over the properties of an object.\" ---\nconst forInParagraphXml =\n  '<w:p w14:paraId=\"606127DB\" w14:textId=\"1704774E\" w:rsidR=\"00B648DA\" w:rsidRDefault=\"00B648DA\" w:rsidP=\"007D4C00\">' +\n  '<w:pPr><w:pStyle w:val=\"NoSpacing\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"9\"/></w:numPr></w:pPr>' +\n  '<w:r><w:t xml:space=\"preserve\">Use </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:proofErr w:type=\"gramStart\"/>' +\n  '<w:r><w:t>for..</w:t></w:r>' +\n  '<w:proofErr w:type=\"gramEnd\"/>' +\n  '<w:r><w:t>in</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> loop to iterate over the properties of an object.</w:t></w:r>' +\n  '</w:p>';\n\nconst forInParagraph = paragraphs.items[23];\nforInParagraph.getRange().insertOoxml(flatOpcPackage(forInParagraphXml), Word.InsertLocation.replace);\nawait context.sync();\n\n// --- Next paragraph in the same bulleted list: was empty, now gets the\n//     \"continue\" keyword note with the relocated _GoBack bookmark. ---\nconst continueParagraphXml =\n  '<w:p w14:paraId=\"76AC9543\" w14:textId=\"77777777\" w:rsidR=\"001A4224\" w:rsidRDefault=\"001A4224\" w:rsidP=\"007D4C00\">' +\n  '<w:pPr><w:pStyle w:val=\"NoSpacing\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"9\"/></w:numPr></w:pPr>' +\n  '<w:r><w:t>The continue key word jumps you to next iteration. This is old legacy JavaScript</w:t></w:r>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> and do not recommend using.</w:t></w:r>' +\n  '</w:p>';\n\n// Re-fetch paragraphs: the previous insertOoxml replaced paragraph 23's\n// content, but the collection itself (paragraph count/order) is unchanged,\n// so the original items/index are still valid to reuse.\nconst continueParagraph = paragraphs.items[24];\ncontinueParagraph.getRange().insertOoxml(flatOpcPackage(continueParagraphXml), Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# PowerShell/Word-COM edit script\n# Target change (per the commit \"9/21/18 - up to exercises\"):\n#   1) In the \"Use for..in loop ...\" bullet, split the single \"for..in\"\n#      run into \"for..\" / \"in\" and wrap \"for..\" in its own\n#      gramStart/gramEnd proofErr pair (matching the existing\n#      \"for..of\" / \"if..else\" bullets elsewhere in the doc), and drop\n#      the stray _GoBack bookmark that Word had left at the end of\n#      that paragraph.\n#   2) The next (previously empty) bullet under the same list gets new\n#      text about the `continue` keyword, with the _GoBack bookmark\n#      relocated into the middle of that new sentence.\n#\n# We use Range.InsertXML with a literal <w:p> fragment so we can\n# control the run/proofErr/bookmark structure exactly, while keeping\n# each paragraph's own pPr (style + numbering) and w14 identity\n# attributes intact.\n\n$d = $word.ActiveDocument\n\n$wNs  = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"'\n$w14Ns = 'xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\"'\n\n# --- Paragraph: \"Use for..in loop to iterate over the properties of an object.\" ---\n$forInPara = $d.Paragraphs(24)\n\n$forInXml = @\"\n<w:p $wNs $w14Ns w14:paraId=\"606127DB\" w14:textId=\"1704774E\" w:rsidR=\"00B648DA\" w:rsidRDefault=\"00B648DA\" w:rsidP=\"007D4C00\"><w:pPr><w:pStyle w:val=\"NoSpacing\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"9\"/></w:numPr></w:pPr><w:r><w:t xml:space=\"preserve\">Use </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:proofErr w:type=\"gramStart\"/><w:r><w:t>for..</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t>in</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> loop to iterate over the properties of an object.</w:t></w:r></w:p>\n\"@\n\n$forInPara.Range.InsertXML($forInXml)\n\n# --- Next paragraph in the same bulleted list: was empty, now gets the\n#     \"continue\" keyword note with the relocated _GoBack bookmark. ---\n$continuePara = $d.Paragraphs(25)\n\n$continueXml = @\"\n<w:p $wNs $w14Ns w14:paraId=\"76AC9543\" w14:textId=\"77777777\" w:rsidR=\"001A4224\" w:rsidRDefault=\"001A4224\" w:rsidP=\"007D4C00\"><w:pPr><w:pStyle w:val=\"NoSpacing\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"9\"/></w:numPr></w:pPr><w:r><w:t>The continue key word jumps you to next iteration. This is old legacy JavaScript</w:t></w:r><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/><w:r><w:t xml:space=\"preserve\"> and do not recommend using.</w:t></w:r></w:p>\n\"@\n\n$continuePara.Range.InsertXML($continueXml)\n"}
